$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.395.99"
$ws.Range("E2").Value = "  +1.47%  "

$ws.Range("D3").Value = "2.275.96"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.495"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").Value = "2.628.70"
$ws.Range("E14").Value = "  +0.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("D16").Value = "2.273.74"
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.794"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.21%  "

$ws.Range("D18").Value = "42.314.38"
$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").Value = "0.0₃0912"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.81%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0744"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.48%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.106"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("E40").Value = "  -1.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.44%  "

$ws.Range("E42").Value = "  +13.95%  "

$ws.Range("D43").Value = "1.997.89"
$ws.Range("E43").Value = "  -1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0286"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.01"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("E46").Value = "  +3.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "92.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.33%  "
